$ErrorActionPreference = "Stop"
$d = $word.ActiveDocument

function Replace-ParagraphXml($SearchText, $NewParaXml) {
    $f = $d.Content.Find
    $f.ClearFormatting()
    $found = $f.Execute($SearchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Find failed for: $SearchText"
    }
    $r = $f.Parent
    $p = $r.Paragraphs(1)
    $pr = $p.Range
    $startPos = $pr.Start
    $endPos = $pr.End - 1
    $full = $d.Range($startPos, $endPos)
    $full.Text = ""
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' + '<w:body>' + $NewParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $full.InsertXML($pkg)
}


Replace-ParagraphXml 'Documentatie' '<w:p w14:paraId="090BFEA7" w14:textId="7A3A2302" w:rsidR="009F0FD6" w:rsidRDefault="00F62540"><w:pPr><w:rPr><w:lang w:val="ro-RO"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00F62540"><w:rPr><w:b/><w:bCs/><w:lang w:val="ro-RO"/></w:rPr><w:t>Documentatie</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006209AF" w:rsidRPr="00CA6712"><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:br/></w:r><w:r w:rsidR="006209AF" w:rsidRPr="00CA6712"><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:br/></w:r><w:hyperlink r:id="rId6" w:history="1"><w:r w:rsidR="00F171C9" w:rsidRPr="00CA6712"><w:rPr><w:rStyle w:val="Hyperlink"/><w:lang w:val="ro-RO"/></w:rPr><w:t>https://glenn-viroux.medium.com/creating-a-music-genre-classifier-using-a-convolutional-neural-network-548d06658cee</w:t></w:r></w:hyperlink><w:r w:rsidR="00EB2076"><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>'

Replace-ParagraphXml 'Se determină valoarea K prin trial-and-error, ca să obținem cea mai bună performanță' '<w:p w14:paraId="7BCBAB5D" w14:textId="456DAB29" w:rsidR="00CA6712" w:rsidRPr="00CA6712" w:rsidRDefault="00CA6712" w:rsidP="00CA6712"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:lang w:val="ro-RO"/></w:rPr></w:pPr><w:r w:rsidRPr="00CA6712"><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>Se determină valoarea K prin trial-and-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>error</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>, ca să obținem cea mai bună performanță</w:t></w:r></w:p>'

Replace-ParagraphXml 'Support Vector Machine' '<w:p w14:paraId="35699DA4" w14:textId="2831F68E" w:rsidR="00CA6712" w:rsidRDefault="00CA6712" w:rsidP="00CA6712"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:lang w:val="ro-RO"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>Support</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve"> Vector </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>Machine</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'

Replace-ParagraphXml 'Clasificare lineară, datele trebuie analizate pe bază de o formulă specificată utilizând un kernel anume, ales pentru cazul de față' '<w:p w14:paraId="2EE0010D" w14:textId="25898BEC" w:rsidR="00CA6712" w:rsidRPr="00CA6712" w:rsidRDefault="00CA6712" w:rsidP="00CA6712"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:rPr><w:lang w:val="ro-RO"/></w:rPr></w:pPr><w:r w:rsidRPr="00CA6712"><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve">Clasificare lineară, datele trebuie analizate pe bază de o formulă specificată utilizând un </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>kernel</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve"> anume, ales pentru cazul de față</w:t></w:r></w:p>'

Replace-ParagraphXml 'FF Neural Network' '<w:p w14:paraId="1EF07216" w14:textId="0FF9BB2F" w:rsidR="00F62540" w:rsidRPr="00F62540" w:rsidRDefault="00CA6712" w:rsidP="00F62540"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:lang w:val="ro-RO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve">FF Neural </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>Network</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'

Replace-ParagraphXml 'Fully connected neural network, funcție de activare cea mai folosită: ReLu' '<w:p w14:paraId="0ADF7D71" w14:textId="5F611598" w:rsidR="002E0690" w:rsidRDefault="002E0690" w:rsidP="002E0690"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr><w:rPr><w:lang w:val="ro-RO"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>Fully</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>connected</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve"> neural </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>network</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve">, funcție de activare cea mai folosită: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>ReLu</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'

Replace-ParagraphXml 'Trebuie folosite metode the regularizare (dropout, L2 etc.) pentru rezultate mai bune' '<w:p w14:paraId="54B69E18" w14:textId="0BC54867" w:rsidR="002773B1" w:rsidRDefault="002773B1" w:rsidP="002E0690"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr><w:rPr><w:lang w:val="ro-RO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>Trebuie folosite metode the regularizare (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>dropout</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>, L2 etc.) pentru rezultate mai bune</w:t></w:r></w:p>'

Replace-ParagraphXml 'Exemplu specific: MLPClassifier, rezultate din (4): acuratețe 89%' '<w:p w14:paraId="5139EE26" w14:textId="6C1387D8" w:rsidR="00F62540" w:rsidRDefault="00F62540" w:rsidP="002E0690"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr><w:rPr><w:lang w:val="ro-RO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve">Exemplu specific: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>MLPClassifier</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>, rezultate din (4): acuratețe 89%</w:t></w:r></w:p>'

Replace-ParagraphXml 'Random Forest' '<w:p w14:paraId="6A8CF4A4" w14:textId="55F2B4BD" w:rsidR="00416CB9" w:rsidRDefault="00416CB9" w:rsidP="00416CB9"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:lang w:val="ro-RO"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>Random</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve"> Forest</w:t></w:r></w:p>'

Replace-ParagraphXml 'Fiecare strat din rețea are substrat de pooling, regularizare, iar în a doua fază sunt câteva straturi de fully connected, urmând straturi de finalizare (softmax output, cross entropy)' '<w:p w14:paraId="26965F16" w14:textId="7D693D13" w:rsidR="002773B1" w:rsidRDefault="002773B1" w:rsidP="002773B1"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr><w:rPr><w:lang w:val="ro-RO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve">Fiecare strat din rețea are substrat de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>pooling</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve">, regularizare, iar în a doua fază sunt câteva straturi de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>fully</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>connected</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>, urmând straturi de finalizare (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>softmax</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve"> output, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>cross</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>entropy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>)</w:t></w:r></w:p>'

Replace-ParagraphXml 'Dataset: ' '<w:p w14:paraId="45DDC7A8" w14:textId="7C4213AD" w:rsidR="007E2F10" w:rsidRPr="007E2F10" w:rsidRDefault="00325239" w:rsidP="00185F95"><w:pPr><w:rPr><w:lang w:val="ro-RO"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>Dataset</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r w:rsidRPr="00325239"><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>GTZAN</w:t></w:r><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>'

Replace-ParagraphXml 'Feature extraction' '<w:p w14:paraId="04C5D6EB" w14:textId="22F94764" w:rsidR="00DC1278" w:rsidRDefault="00885123" w:rsidP="00DC1278"><w:pPr><w:rPr><w:lang w:val="ro-RO"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>Feature</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>extraction</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'

Replace-ParagraphXml 'În general este folosit un color-map pentru reprezentarea grafică' '<w:p w14:paraId="63F587C9" w14:textId="021C2BD0" w:rsidR="00885123" w:rsidRDefault="00885123" w:rsidP="00885123"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="8"/></w:numPr><w:rPr><w:lang w:val="ro-RO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>În general este folosit un color-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>map</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve"> pentru reprezentarea grafică</w:t></w:r></w:p>'

Replace-ParagraphXml ' Spectrograma Mel este ' '<w:p w14:paraId="1F05029A" w14:textId="1C2693DF" w:rsidR="00885123" w:rsidRDefault="00A17054" w:rsidP="00885123"><w:pPr><w:rPr><w:lang w:val="ro-RO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>Pentru clasificare genuri: conform (</w:t></w:r><w:r w:rsidR="00CE34F4"><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>) și (6) =</w:t></w:r><w:r><w:t>&gt;</w:t></w:r><w:r w:rsidRPr="003D2F29"><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Spectrograma</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> Mel </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>este</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="003D2F29"><w:rPr><w:u w:val="single"/><w:lang w:val="ro-RO"/></w:rPr><w:t>mai eficientă</w:t></w:r></w:p>'

Replace-ParagraphXml 'GTZAN dataset, Spectograma Mel, CNN+SVM' '<w:p w14:paraId="2E5C86AB" w14:textId="0ED5D86D" w:rsidR="00BB2AF2" w:rsidRDefault="00BB2AF2" w:rsidP="00885123"><w:pPr><w:rPr><w:lang w:val="ro-RO"/></w:rPr></w:pPr><w:r w:rsidRPr="00E95BAC"><w:rPr><w:b/><w:bCs/><w:lang w:val="ro-RO"/></w:rPr><w:t>Metode alese</w:t></w:r><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:br/><w:t xml:space="preserve">GTZAN </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>dataset</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>Spectograma</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve"> Mel, CNN+SVM</w:t></w:r></w:p>'


Write-Output "All replacements applied."
